$d = $word.ActiveDocument

$replacements = @(
    @("65÷6=", "93÷7="),
    @("52÷5=", "86÷3="),
    @("52÷6=", "43÷9="),
    @("49÷2=", "98÷3="),
    @("44÷7=", "78÷3="),
    @("62÷8=", "78÷8="),
    @("47÷7=", "24÷8="),
    @("74÷3=", "41÷9="),
    @("45÷9=", "13÷3="),
    @("28÷4=", "52÷2="),
    @("73÷3=", "64÷3="),
    @("10÷6=", "41÷9="),
    @("42÷2=", "24÷8="),
    @("17÷4=", "76÷8="),
    @("95÷6=", "22÷7="),
    @("41÷7=", "73÷4="),
    @("75÷2=", "86÷3="),
    @("46÷2=", "64÷5="),
    @("48÷3=", "50÷8="),
    @("22÷5=", "49÷4="),
    @("73÷8=", "31÷9="),
    @("47÷3=", "70÷5="),
    @("19÷9=", "14÷8="),
    @("36÷9=", "85÷3="),
    @("40÷2=", "83÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
